$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Recepcionista"
$ws.Range("B1").Value = "18_12_2023"
$ws.Range("C1").Value = "05_01_2024"
$ws.Range("D1").Value = "15_01_2024"
$ws.Range("E1").Value = "21_01_2024"
$ws.Range("F1").Value = "28_01_2024"
$ws.Range("G1").Value = "07_02_2024"
$ws.Range("H1").Value = "11_02_2024"
$ws.Range("I1").Value = "18_02_2024"

# Row 2 - Alejandro
$ws.Range("A2").Value = "Alejandro"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 13
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 1

# Row 3 - Camila
$ws.Range("A3").Value = "Camila"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 14
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 5

# Row 4 - Betty
$ws.Range("A4").Value = "Betty"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 3

# Row 5 - Felipe
$ws.Range("A5").Value = "Felipe"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 4

# Row 6 - Constanza
$ws.Range("A6").Value = "Constanza"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 2

$ws.Range("B6").Font.Underline = $true
$ws.Range("D6").Font.Underline = $true
$ws.Range("F6").Font.Underline = $true
$ws.Range("H6").Font.Underline = $true

# Row 8 - stray formatted (underlined) empty cell at C8
$ws.Range("C8").Font.Underline = $true

# Column widths: column A already carries the target width from the
# original file, so leave it untouched and only size B (to match A) and C.
$ws.Columns("B").ColumnWidth = 11.5
$ws.Columns("C").ColumnWidth = 11.71

$ws.Range("I4").Select() | Out-Null
